# Duplicate the question block (rows 2:4) twice, into rows 5:7 and 8:10,
# matching the exam-search / countdown-clock fix that grew the sheet from
# one question (rows 1-4) to three stacked questions (rows 1-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 5:7 = copy of rows 2:4 (formats, then values) ---
$ws.Range("A2:J4").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:J4").Copy()
$ws.Range("A5").PasteSpecial(-4104)   # xlPasteValues

# --- Block 2: rows 8:10 = copy of rows 2:4 (formats, then values) ---
$ws.Range("A2:J4").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:J4").Copy()
$ws.Range("A8").PasteSpecial(-4104)   # xlPasteValues

# Row heights for the newly created rows match the existing question rows.
$ws.Rows.Item(5).RowHeight = 49.2
$ws.Rows.Item(6).RowHeight = 49.2
$ws.Rows.Item(7).RowHeight = 49.2
$ws.Rows.Item(8).RowHeight = 49.2
$ws.Rows.Item(9).RowHeight = 49.2
$ws.Rows.Item(10).RowHeight = 49.2

# Extend the "0/1" list validation from the original 2:4 rows down to 2:10
# for every answer's correctness column (D, F, H, J). Also re-create the
# D1 custom validation afterwards so it keeps sorting after the list rules
# (matching the original list-then-custom ordering).
$ws.Range("D2:D4").Validation.Delete()
$ws.Range("F2:F4").Validation.Delete()
$ws.Range("H2:H4").Validation.Delete()
$ws.Range("J2:J4").Validation.Delete()
$ws.Range("D1").Validation.Delete()
$ws.Range("D2:D10").Validation.Add(3, 1, 1, '"0,1"')
$ws.Range("F2:F10").Validation.Add(3, 1, 1, '"0,1"')
$ws.Range("H2:H10").Validation.Add(3, 1, 1, '"0,1"')
$ws.Range("J2:J10").Validation.Add(3, 1, 1, '"0,1"')
$ws.Range("D1").Validation.Add(7, 1, 1, "1")

# Restore the usual view/selection state (new last three rows selected).
$ws.Range("A8:A10").EntireRow.Select()
